# The deck had the "Code Review" slide (position 7) immediately followed by
# the "Take away" slide (position 8). The author reordered them so that
# "Take away" now comes before "Code Review" (i.e. the two slides swap
# places; no textual content on either slide changes).
$p = $ppt.ActivePresentation

# Locate the two slides by their title text so the script is resilient to
# being run against a deck where they aren't exactly at indices 7/8.
$takeAway = $null
$codeReview = $null
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    $title = ""
    try { $title = $slide.Shapes.Title.TextFrame.TextRange.Text } catch { $title = "" }
    if ($title -like "Take*away*") { $takeAway = $slide }
    if ($title -like "Code*Review*") { $codeReview = $slide }
}

if ($takeAway -ne $null -and $codeReview -ne $null -and $takeAway.SlideIndex -gt $codeReview.SlideIndex) {
    # Move "Take away" to just before "Code Review".
    $takeAway.MoveTo($codeReview.SlideIndex)
}
